$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.076.98'
$ws.Range('E2').Value = '  +0.42%  '
$ws.Range('D3').Value = '2.222.31'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '291.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '87.80'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.89%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.512'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.473'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '30.42'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0780'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E12').Value = '  +3.28%  '
$ws.Range('E13').Value = '  +1.48%  '
$ws.Range('D14').Value = '2.568.92'
$ws.Range('E14').Value = '  -0.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '13.97'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.11%  '
$ws.Range('D16').Value = '2.235.35'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.728'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').Value = '40.036.98'
$ws.Range('E18').Value = '  +0.43%  '
$ws.Range('D19').Value = '0.0₃0885'
$ws.Range('E19').Value = '  -0.91%  '
$ws.Range('E20').Value = '  +7.30%  '
$ws.Range('E21').Value = '  +0.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.65'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.79'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.74%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.46'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.53%  '
$ws.Range('E26').Value = '  -0.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.71'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.19'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '156.71'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.55%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '31.80'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.95'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0718'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.21%  '
$ws.Range('E35').Value = '  -0.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.87'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.81%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '15.64'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0982'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.66%  '
$ws.Range('E40').Value = '  +2.00%  '
$ws.Range('D41').Value = '2.117.71'
$ws.Range('E41').Value = '  +8.14%  '
$ws.Range('E42').Value = '  +2.40%  '
$ws.Range('E43').Value = '  -2.09%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '18.02'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +11.73%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0268'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.85%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.95'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.66'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.63%  '
$ws.Range('D48').Value = '2.437.19'
$ws.Range('E48').Value = '  -0.54%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.45'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '69.57'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.86%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.10'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.06%  '
